# PlayerPerformance_5955.xlsx — add "ODI Batting Extra" / "ODI Bowling Extra"
# sheets and tidy up stray empty cells in "ODI Batting" (commit:
# "[AFG] added final excel sheets for Afghanistan").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "ODI Batting": clear the handful of empty placeholder cells in column B
#    (rows where INNING_NUMBER was never populated).
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
foreach ($r in 2,5,6,7,8,9,11) {
    $batting.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 2. Add the two new sheets at the end of the workbook, in order.
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add($null, $afterSheet)
$battingExtra.Name = "ODI Batting Extra"

$afterSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $afterSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header style used throughout the workbook (bold, thin border, centered).
$headerStyleSource = $wb.Worksheets.Item("Player Info").Range("A1")

# ---------------------------------------------------------------------------
# 3. "ODI Batting Extra" content.
# ---------------------------------------------------------------------------
$battingHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($i = 0; $i -lt $battingHeaders.Length; $i++) {
    $battingExtra.Cells.Item(1, $i + 1).Value = $battingHeaders[$i]
}
$headerStyleSource.Copy()
$battingExtra.Range("A1:F1").PasteSpecial(-4122)

# Columns: MATCH_CODE(text), BATTING_POSITION(number), NUM_4(text), NUM_6(text),
# PERCENT_RUNS_OF_TOTAL(text), MAN_OF_MATCH(text). A leading apostrophe keeps
# numeric-looking strings stored as text instead of being coerced to numbers.
$battingRows = @(
    @("'4530", 11,   $null,  $null,  $null,     "NO"),
    @("'4537", 11,   "'0",   "'0",   $null,     "NO"),
    @("'4538", $null,$null,  $null,  $null,     "NO"),
    @("'4539", 10,   $null,  $null,  $null,     "NO"),
    @("'4582", 10,   $null,  $null,  $null,     "NO"),
    @("'4585", 10,   $null,  $null,  $null,     "NO"),
    @("'4588", $null,$null,  $null,  $null,     "NO"),
    @("'4671", 10,   "'0",   "'0",   "'0.88%",  "NO"),
    @("'4674", $null,$null,  $null,  $null,     "NO"),
    @("'4675", $null,$null,  $null,  $null,     $null)
)

for ($i = 0; $i -lt $battingRows.Length; $i++) {
    $row = $i + 2
    $data = $battingRows[$i]

    $battingExtra.Cells.Item($row, 1).Value = $data[0]

    if ($null -ne $data[1]) {
        $battingExtra.Cells.Item($row, 2).Value = $data[1]
    }
    if ($null -ne $data[2]) {
        $battingExtra.Cells.Item($row, 3).Value = $data[2]
    }
    if ($null -ne $data[3]) {
        $battingExtra.Cells.Item($row, 4).Value = $data[3]
    }
    if ($null -ne $data[4]) {
        $battingExtra.Cells.Item($row, 5).Value = $data[4]
    }
    if ($null -ne $data[5]) {
        $battingExtra.Cells.Item($row, 6).Value = $data[5]
    }
}

# ---------------------------------------------------------------------------
# 4. "ODI Bowling Extra" content.
# ---------------------------------------------------------------------------
$bowlingHeaders = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($i = 0; $i -lt $bowlingHeaders.Length; $i++) {
    $bowlingExtra.Cells.Item(1, $i + 1).Value = $bowlingHeaders[$i]
}
$headerStyleSource.Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

$bowlingRows = @(
    @("'4530", "'1", "'10.00%"),
    @("'4537", "'1", "'40.00%"),
    @("'4538", $null, $null),
    @("'4539", "'0", "'10.00%"),
    @("'4582", "'1", "'20.00%"),
    @("'4585", "'1", "'20.00%"),
    @("'4588", "'0", "'20.00%"),
    @("'4671", $null, $null),
    @("'4674", "'0", $null),
    @("'4675", $null, $null)
)

for ($i = 0; $i -lt $bowlingRows.Length; $i++) {
    $row = $i + 2
    $data = $bowlingRows[$i]

    $bowlingExtra.Cells.Item($row, 1).Value = $data[0]

    if ($null -ne $data[1]) {
        $bowlingExtra.Cells.Item($row, 2).Value = $data[1]
    }
    if ($null -ne $data[2]) {
        $bowlingExtra.Cells.Item($row, 3).Value = $data[2]
    }
}

# Keep the originally-active sheet selected (workbook.xml's activeTab="0").
$wb.Worksheets.Item("Player Info").Activate()
